$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E ("% of total people"), shifting the old "average debt
# per person" column (and its data) to column F.
$ws.Range("E1").EntireColumn.Insert()

$ws.Range("E1").Value = "% of total people"

# Reorder the data rows (2-8) sorted by total debt (column C) descending,
# and fill in the new "% of total people" column.
$names   = @("merchant", "shoemaker", "executor", "frederick company", "shopkeeper", "administrator", "esq")
$totals  = @(78510.94, 1237.07, 303.77, 206.67, 33.33, 22.5, 17.665733)
$counts  = @(47, 2, 3, 1, 1, 2, 1)
$percent = @(82.45614035087719, 3.508771929824561, 5.263157894736842, 1.754385964912281, 1.754385964912281, 3.508771929824561, 1.754385964912281)
$avgs    = @(1670.445531914894, 618.535, 101.2566666666667, 206.67, 33.33, 11.25, 17.665733)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $totals[$i]
    $ws.Cells.Item($r, 4).Value = $counts[$i]
    $ws.Cells.Item($r, 5).Value = $percent[$i]
    $ws.Cells.Item($r, 6).Value = $avgs[$i]
}
